$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# A new "Kelp" flag column is being inserted at column G (a formula
# =OR(Macro=1,Nereo=1)). The handful of free-text annotations that used
# to live in column G (and the one in H8) shift one column to the right,
# to make room. We move those first (copying formatting along with the
# text), then overwrite column G with the new formula column.
# ------------------------------------------------------------------

# Row 8: ", " note moves from H8 to I8 (keep its distinct style).
$ws.Range("H8").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$ws.Range("I8").Value = ", "
$ws.Range("H8").ClearContents()

# Row 5: G5 was an empty, specially-styled cell; that formatting moves to H5.
$ws.Range("G5").Copy()
$ws.Range("H5").PasteSpecial(-4122)

# Rows 65 & 67: "Line connecting ..." notes shift from G to H verbatim.
$ws.Range("H65").Value = "Line connecting nanat surveys"
$ws.Range("G65").ClearContents()

$ws.Range("H67").Value = "Line connecting KB4 sites"
$ws.Range("G67").ClearContents()

# Rows 3, 4, 19, 43: "Coordinates fixed" notes are simply dropped (no
# longer needed) - nothing moves to H for these, column G below gets the
# new formula instead.
$ws.Range("G3").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("G19").ClearContents()
$ws.Range("G43").ClearContents()

# ------------------------------------------------------------------
# New column G: header "Kelp" + boolean formula flagging rows present in
# either the Macro or Nereo column.
# ------------------------------------------------------------------
$ws.Range("G1").Value = "Kelp"

# Row 40: note text moves to H40, with updated wording (replicate note).
$ws.Range("H40").Value = 'Note to throw this data point should not be used due to replicate with "27"'
$ws.Range("G40").ClearContents()

$ws.Range("G2").Formula = "=OR(E2=1,F2=1)"
$ws.Range("G3:G66").Formula = "=OR(E3=1,F3=1)"
$ws.Range("G67:G68").Formula = "=OR(E67=1,F67=1)"
# G5 inherited the old s="2" styling when the formula filled across the
# range; reset it back to the plain/default style used by the rest of
# the new column.
$ws.Range("G5").Style = "Normal"

# ------------------------------------------------------------------
# View/selection bookkeeping to mirror the saved workbook state.
# ------------------------------------------------------------------
$ws.Range("H41").Select() | Out-Null
